$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Garsoniera str.Parang, et.3, mobilat, utilat", "180 EUR / lună", "https://www.imobiliare.ro/inchirieri-garsoniere/cluj-napoca/manastur/garsoniera-de-inchiriat-XB7200028"),
    @("Inchiriere Apartament 1 camera semidecomandat, 25 mp, Etajul 4 din 4", "200 EUR / lună", "https://www.imobiliare.ro/inchirieri-garsoniere/cluj-napoca/manastur/garsoniera-de-inchiriat-X01V104V9"),
    @("Garsoniera, Manastur", "200 EUR / lună", "https://www.imobiliare.ro/inchirieri-garsoniere/cluj-napoca/manastur/garsoniera-de-inchiriat-X9PL1035U"),
    @("Apartament 1 camera str.Bucegi, zona McDonald;s, decomandat", "200 EUR / lună", "https://www.imobiliare.ro/inchirieri-garsoniere/cluj-napoca/manastur/garsoniera-de-inchiriat-XARU0004J"),
    @("Garsoniera in zona Piata Flora", "230 EUR / lună", "https://www.imobiliare.ro/inchirieri-garsoniere/cluj-napoca/manastur/garsoniera-de-inchiriat-X8M2103DP")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("A1:A5").Copy()
$ws.Range("B1:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
